$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update percentages for remaining rows (keep as literal text, not numeric percent)
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "50.00%"
$ws.Range("B2").Style = "Normal"

$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "50.00%"
$ws.Range("B3").Style = "Normal"

# Remove the "pendiente" row entirely (row 4)
$ws.Rows.Item(4).Delete()
